# PostgreSQL-persons.xlsx update:
#   - "create-accounts" sheet gains a new "is_active boolean NOT NULL DEFAULT true"
#     column definition row, inserted just above the closing ");" row.
#   - "insert_persons" sheet's three "INSERT INTO private.accounts" rows get a
#     trailing ", true" value added for the new column.
#   - The active/selected worksheet moves from "create-accounts" to "insert_persons".

$wb = $excel.ActiveWorkbook

# ---- Sheet: create-accounts ---------------------------------------------
$ws2 = $wb.Worksheets.Item("create-accounts")

# Insert a new blank row above row 6 (the "');'" closing row), pushing it to row 7.
$ws2.Range("A6").EntireRow.Insert()

# Row 5 ("updated_at ... NOT NULL") is no longer the last column definition,
# so it needs a trailing comma.
$ws2.Range("C5").Value = "NOT NULL,"

# Fill in the newly inserted row 6 with the new column definition.
$ws2.Range("A6").Value = "is_active"
$ws2.Range("B6").Value = "boolean"
$ws2.Range("C6").Value = "NOT NULL DEFAULT true"

# Matches the saved selection on this sheet after the edit.
[void]$ws2.Range("C6").Select()

# ---- Sheet: insert_persons -----------------------------------------------
$ws3 = $wb.Worksheets.Item("insert_persons")

foreach ($r in 5..7) {
    $cell = $ws3.Cells.Item($r, 1)
    $cell.Value = $cell.Text.Replace("'2019-08-02')", "'2019-08-02', true)")
}

# The workbook now opens with "insert_persons" as the active tab.
$ws3.Activate()
